# The "March" attendance sheet was missing the column for 2023-03-13 —
# `add_column` had an off-by-one bug that returned the wrong column
# number, so the day's data never landed in the sheet. Add it as the
# next column (D): a date header plus the student's in/out times.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("March")

# D2 (the in/out-time text) isn't date-like, so a direct Value assign
# stores it as plain text exactly as typed.
$ws.Cells.Item(2, 4).Value = "In-time: 00:04:45 `nOut-time: 00:04:49"

# D1's header text ("2023-03-13") looks like a date, and Excel would
# normally reinterpret a directly-assigned date-like string as a real
# date serial. Stage it as explicit text in a scratch cell, then copy/
# paste-values it into place so the destination keeps plain, unstyled
# text - just like the rest of the header row.
$scratch = $ws.Cells.Item(5, 10)
$scratch.Value = "'2023-03-13"
$scratch.Copy()
$ws.Cells.Item(1, 4).PasteSpecial(-4163)
$scratch.Clear()

# The header row grows slightly taller to fit the new column.
$ws.Rows.Item(1).RowHeight = 19.5
